$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "barrios.zip"
$ws.Range("B12").Value = "https://recursos-data.buenosaires.gob.ar/ckan2/barrios/barrios.zip"
$ws.Range("E12").Value = "http://data.buenosaires.gob.ar/dataset/barrios"

$excel.ActiveWindow.ScrollRow = 6
$ws.Range("E12").Select()
